# Update the "Förändrad" (Changed) date column C for all data rows
# (rows 2-141) from serial date 45188 to 45189 (i.e. +1 day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 141) { $lastRow = 141 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
